# Auto-generated Excel COM-interop edit script
# Applies the cell-level changes described by the commit diff to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '245.03'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '4.999'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '6.576'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8114'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.8435'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1338'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06943'
$ws.Range('E11').Value = '10MandalaExchangeTokenMDX'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.02839'
$ws.Range('E12').Value = '11BitrueCoinBTR'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.09403'
$ws.Range('E13').Value = '12BitMartTokenBMX'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.001516'
$ws.Range('E14').Value = '13BitForexTokenBF'
$ws.Range('B15').Value = 'One'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0005959'
$ws.Range('E15').Value = '14OneONEWorstin24h'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.006087'
$ws.Range('E16').Value = '15TigerCashTCH'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.499'
$ws.Range('E17').Value = '16LEOLEO'
$ws.Range('B18').Value = 'BTSEToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.092'
$ws.Range('E18').Value = '17BTSETokenBTSE'
$ws.Range('B19').Value = 'BitpandaEcosystemToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.3166'
$ws.Range('E19').Value = '18BitpandaEcosystemTokenBEST'
$ws.Range('B20').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C20').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.03288'
$ws.Range('E20').Value = '19LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1319'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.736'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.004524'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.00009695'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0001939'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03661'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1356'
$ws.Range('B42').Value = 'KickToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.006237'
$ws.Range('E42').Value = '41KickTokenKICK'
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002722'
$ws.Range('E43').Value = '42CEJICEJI'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.008089'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005294'
$ws.Range('E47').Value = '46CoinbaseStockTokenCOIN'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.002040'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002099'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0001999'

Write-Host "Applied 67 cell updates."
